$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "A1" = -0.23237242672377789
    "B1" = 0.231887096834555
    "A2" = -0.16780440163650834
    "B2" = 0.16662856493719591
    "A3" = -0.097241244567197072
    "B3" = 0.097015783317774051
    "A4" = -0.089015783405843152
    "B4" = 0.088600379344919489
    "A5" = -0.085600379394421111
    "B5" = 0.084192133621818321
    "A6" = 0.0150666658501315
    "B6" = -0.015179199839865376
    "A7" = 0.025179199719811418
    "B7" = -0.025202838728056953
    "A8" = 0.035202838610819409
    "B8" = -0.035268706372891589
    "A9" = 0.037268706330739754
    "B9" = -0.037337613742499443
    "A10" = 0.039337613706980079
    "B10" = -0.039342119723935198
    "A11" = 0.042342119680442991
    "B11" = -0.042356150160621553
    "A12" = 0.045856150114288674
    "B12" = -0.046029443387143765
    "A13" = -0.012877011167848451
    "B13" = 0.012868333237761043
    "A14" = -0.0048683333114558636
    "B14" = 0.0048667450913706389
    "A15" = -0.0038667451050509172
    "B15" = 0.0038611320388497461
    "A16" = 0.006819691479500456
    "B16" = -0.0072225174061606445
    "A17" = 0.0092225173852797937
    "B17" = -0.0093853686574361106
    "A18" = -0.058082695252128502
    "B18" = 0.057933312867692877
    "A19" = -0.053933312905632746
    "B19" = 0.052803037603866176
    "A20" = -0.048803037655085646
    "B20" = 0.048471538890588306
    "A21" = -0.0040057842098297769
    "B21" = 0.0039999999454440882
    "A22" = -0.065387320419285189
    "B22" = 0.065080963888775045
    "A23" = -0.060080963950879251
    "B23" = 0.05948578184542086
    "A24" = -0.039485782045134421
    "B24" = 0.039279031126285169
    "A25" = -0.097258799524153616
    "B25" = 0.09713475488813117
    "A26" = -0.09463475494749396
    "B26" = 0.094474198450116731
    "A27" = -0.091974198513739669
    "B27" = 0.09101774312951294
    "A28" = -0.089017743207945976
    "B28" = 0.088363051869285947
    "A29" = -0.081363052004462588
    "B29" = 0.081173617285585919
    "A30" = -0.021173617883518681
    "B30" = 0.021022747804721131
    "A31" = -0.014022747950532377
    "B31" = 0.014000894751239557
    "A32" = -0.0040008949234771052
    "B32" = 0.003999999879280125
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
